# Auto-generated edit script: update computed market-price columns
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the 8 crafting-job sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 58.23077
$ws.Cells.Item(11, 9).Value = 58.23077
$ws.Cells.Item(11, 11).Value = 58.23077
$ws.Cells.Item(11, 13).Value = 81.76922999999999
# Row 33
$ws.Cells.Item(33, 8).Value = 357.30768
$ws.Cells.Item(33, 9).Value = 374.3
$ws.Cells.Item(33, 10).Value = 300.66666
$ws.Cells.Item(33, 11).Value = 374.3
$ws.Cells.Item(33, 12).Value = 300.66666
$ws.Cells.Item(33, 13).Value = -145.3
$ws.Cells.Item(33, 14).Value = -758.66666
# Row 69
$ws.Cells.Item(69, 8).Value = 3000
$ws.Cells.Item(69, 9).Value = 3000
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 9000
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -8126
# Row 72
$ws.Cells.Item(72, 8).Value = 3000
$ws.Cells.Item(72, 9).Value = 3000
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 27000
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -22632
# Row 74
$ws.Cells.Item(74, 8).Value = 3766.3333
$ws.Cells.Item(74, 9).Value = 3286.75
$ws.Cells.Item(74, 11).Value = 3286.75
$ws.Cells.Item(74, 13).Value = -2350.75
# Row 77
$ws.Cells.Item(77, 8).Value = 3766.3333
$ws.Cells.Item(77, 9).Value = 3286.75
$ws.Cells.Item(77, 11).Value = 16433.75
$ws.Cells.Item(77, 13).Value = -11753.75
# Row 112
$ws.Cells.Item(112, 8).Value = 1746.7307
$ws.Cells.Item(112, 9).Value = 800
$ws.Cells.Item(112, 10).Value = 1825.625
$ws.Cells.Item(112, 11).Value = 2400
$ws.Cells.Item(112, 12).Value = 5476.875
$ws.Cells.Item(112, 13).Value = -1292
$ws.Cells.Item(112, 14).Value = -7692.875
# Row 133
$ws.Cells.Item(133, 8).Value = 67635
$ws.Cells.Item(133, 10).Value = 67635
$ws.Cells.Item(133, 12).Value = 67635
$ws.Cells.Item(133, 14).Value = -77755
# Row 137
$ws.Cells.Item(137, 8).Value = 1479.7715
$ws.Cells.Item(137, 9).Value = 1284.6666
$ws.Cells.Item(137, 11).Value = 3853.9998
$ws.Cells.Item(137, 13).Value = -1303.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 52
$ws.Cells.Item(52, 8).Value = 27799.75
$ws.Cells.Item(52, 10).Value = 27799.75
$ws.Cells.Item(52, 12).Value = 27799.75
$ws.Cells.Item(52, 14).Value = -28435.75
# Row 74
$ws.Cells.Item(74, 8).Value = 1201.64
$ws.Cells.Item(74, 9).Value = 1066.35
$ws.Cells.Item(74, 10).Value = 1742.8
$ws.Cells.Item(74, 11).Value = 1066.35
$ws.Cells.Item(74, 12).Value = 1742.8
$ws.Cells.Item(74, 13).Value = -192.3499999999999
$ws.Cells.Item(74, 14).Value = -3490.8
# Row 77
$ws.Cells.Item(77, 8).Value = 1201.64
$ws.Cells.Item(77, 9).Value = 1066.35
$ws.Cells.Item(77, 10).Value = 1742.8
$ws.Cells.Item(77, 11).Value = 5331.75
$ws.Cells.Item(77, 12).Value = 8714
$ws.Cells.Item(77, 13).Value = -963.75
$ws.Cells.Item(77, 14).Value = -17450

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 2100.7856
$ws.Cells.Item(134, 9).Value = 1624.9048
$ws.Cells.Item(134, 10).Value = 3528.4285
$ws.Cells.Item(134, 11).Value = 4874.7144
$ws.Cells.Item(134, 12).Value = 10585.2855
$ws.Cells.Item(134, 13).Value = -2339.7144
$ws.Cells.Item(134, 14).Value = -15655.2855

$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Cells.Item(21, 8).Value = 7507.5
$ws.Cells.Item(21, 9).Value = 5000
$ws.Cells.Item(21, 10).Value = 10015
$ws.Cells.Item(21, 11).Value = 5000
$ws.Cells.Item(21, 12).Value = 10015
$ws.Cells.Item(21, 13).Value = -4765
$ws.Cells.Item(21, 14).Value = -10485
# Row 31
$ws.Cells.Item(31, 8).Value = 32262206
$ws.Cells.Item(31, 9).Value = 166673260
$ws.Cells.Item(31, 11).Value = 166673260
$ws.Cells.Item(31, 13).Value = -166672965
# Row 34
$ws.Cells.Item(34, 8).Value = 32262206
$ws.Cells.Item(34, 9).Value = 166673260
$ws.Cells.Item(34, 11).Value = 166673260
$ws.Cells.Item(34, 13).Value = -166673058
# Row 58
$ws.Cells.Item(58, 8).Value = 3898
$ws.Cells.Item(58, 9).Value = 4250
$ws.Cells.Item(58, 10).Value = 3194
$ws.Cells.Item(58, 11).Value = 4250
$ws.Cells.Item(58, 12).Value = 3194
$ws.Cells.Item(58, 13).Value = -4047
$ws.Cells.Item(58, 14).Value = -3600
# Row 99
$ws.Cells.Item(99, 8).Value = 3632.96
$ws.Cells.Item(99, 9).Value = 3683.7368
$ws.Cells.Item(99, 10).Value = 3472.1667
$ws.Cells.Item(99, 11).Value = 3683.7368
$ws.Cells.Item(99, 12).Value = 3472.1667
$ws.Cells.Item(99, 13).Value = -2185.7368
$ws.Cells.Item(99, 14).Value = -6468.1667
# Row 126
$ws.Cells.Item(126, 8).Value = 3632.96
$ws.Cells.Item(126, 9).Value = 3683.7368
$ws.Cells.Item(126, 10).Value = 3472.1667
$ws.Cells.Item(126, 11).Value = 11051.2104
$ws.Cells.Item(126, 12).Value = 10416.5001
$ws.Cells.Item(126, 13).Value = -8581.2104
$ws.Cells.Item(126, 14).Value = -15356.5001
# Row 132
$ws.Cells.Item(132, 8).Value = 2427.9375
$ws.Cells.Item(132, 9).Value = 1579.0834
$ws.Cells.Item(132, 10).Value = 4974.5
$ws.Cells.Item(132, 11).Value = 4737.2502
$ws.Cells.Item(132, 12).Value = 14923.5
$ws.Cells.Item(132, 13).Value = -2207.2502
$ws.Cells.Item(132, 14).Value = -19983.5
# Row 134
$ws.Cells.Item(134, 8).Value = 1724.9286
$ws.Cells.Item(134, 9).Value = 1713.5
$ws.Cells.Item(134, 10).Value = 1753.5
$ws.Cells.Item(134, 11).Value = 5140.5
$ws.Cells.Item(134, 12).Value = 5260.5
$ws.Cells.Item(134, 13).Value = -2605.5
$ws.Cells.Item(134, 14).Value = -10330.5
# Row 136
$ws.Cells.Item(136, 8).Value = 3898
$ws.Cells.Item(136, 9).Value = 4250
$ws.Cells.Item(136, 10).Value = 3194
$ws.Cells.Item(136, 11).Value = 12750
$ws.Cells.Item(136, 12).Value = 9582
$ws.Cells.Item(136, 13).Value = -10200
$ws.Cells.Item(136, 14).Value = -14682

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Cells.Item(81, 8).Value = 4699
$ws.Cells.Item(81, 9).Value = 1498.3334
$ws.Cells.Item(81, 10).Value = 9500
$ws.Cells.Item(81, 11).Value = 4495.0002
$ws.Cells.Item(81, 12).Value = 28500
$ws.Cells.Item(81, 13).Value = -3372.0002
$ws.Cells.Item(81, 14).Value = -30746
# Row 84
$ws.Cells.Item(84, 8).Value = 4699
$ws.Cells.Item(84, 9).Value = 1498.3334
$ws.Cells.Item(84, 10).Value = 9500
$ws.Cells.Item(84, 11).Value = 13485.0006
$ws.Cells.Item(84, 12).Value = 85500
$ws.Cells.Item(84, 13).Value = -7869.000599999999
$ws.Cells.Item(84, 14).Value = -96732
# Row 131
$ws.Cells.Item(131, 8).Value = 32263008
$ws.Cells.Item(131, 9).Value = 11468.889
$ws.Cells.Item(131, 10).Value = 45456820
$ws.Cells.Item(131, 11).Value = 34406.667
$ws.Cells.Item(131, 12).Value = 136370460
$ws.Cells.Item(131, 13).Value = -29366.667
$ws.Cells.Item(131, 14).Value = -136380540

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 12).ClearContents()
$ws.Cells.Item(6, 14).Value = 0
# Row 16
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).ClearContents()
$ws.Cells.Item(16, 14).Value = 0
# Row 132
$ws.Cells.Item(132, 8).Value = 2554
$ws.Cells.Item(132, 9).Value = 1880.6428
$ws.Cells.Item(132, 10).Value = 4439.4
$ws.Cells.Item(132, 11).Value = 5641.928400000001
$ws.Cells.Item(132, 12).Value = 13318.2
$ws.Cells.Item(132, 13).Value = -3111.928400000001
$ws.Cells.Item(132, 14).Value = -18378.2

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Cells.Item(9, 8).Value = 1494.5
$ws.Cells.Item(9, 9).Value = 590
$ws.Cells.Item(9, 10).Value = 2399
$ws.Cells.Item(9, 11).Value = 590
$ws.Cells.Item(9, 12).Value = 2399
$ws.Cells.Item(9, 13).Value = -366
$ws.Cells.Item(9, 14).Value = -2847

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Cells.Item(22, 8).Value = 70007.5
$ws.Cells.Item(22, 10).Value = 90015
$ws.Cells.Item(22, 12).Value = 90015
$ws.Cells.Item(22, 14).Value = -90601
# Row 81
$ws.Cells.Item(81, 8).Value = 73620.07000000001
$ws.Cells.Item(81, 9).Value = 93034.73
$ws.Cells.Item(81, 10).Value = 2433
$ws.Cells.Item(81, 11).Value = 186069.46
$ws.Cells.Item(81, 12).Value = 4866
$ws.Cells.Item(81, 13).Value = -185008.46
$ws.Cells.Item(81, 14).Value = -6988
# Row 84
$ws.Cells.Item(84, 8).Value = 73620.07000000001
$ws.Cells.Item(84, 9).Value = 93034.73
$ws.Cells.Item(84, 10).Value = 2433
$ws.Cells.Item(84, 11).Value = 930347.2999999999
$ws.Cells.Item(84, 12).Value = 24330
$ws.Cells.Item(84, 13).Value = -925043.2999999999
$ws.Cells.Item(84, 14).Value = -34938
# Row 126
$ws.Cells.Item(126, 8).Value = 3703
$ws.Cells.Item(126, 9).Value = 4490.75
$ws.Cells.Item(126, 10).Value = 552
$ws.Cells.Item(126, 11).Value = 13472.25
$ws.Cells.Item(126, 12).Value = 1656
$ws.Cells.Item(126, 13).Value = -11002.25
$ws.Cells.Item(126, 14).Value = -6596
